# Implements the "UPDATE command" change:
#   1. Inserts a new bullet (ilvl=2, numId=2) right before the "DELETE
#      Command" heading, reading:
#         UPDATE primarykey IN tablename SET field TO value
#      with the same run-splitting / proofErr wrapping style used
#      elsewhere in the document (spell-check markers around the
#      made-up identifiers "primarykey" and "tablename").
#   2. Re-splits the trailing " FROM tablename" run of the existing
#      "DELETE primarykey FROM tablename" bullet so "tablename" is
#      wrapped in its own proofErr spellStart/spellEnd run, matching
#      the same convention.

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="20"/></w:rPr>'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $bodyXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1. Insert the new "UPDATE primarykey IN tablename SET field TO value"
#    bullet immediately before the "DELETE Command" paragraph.
# ---------------------------------------------------------------------

$deleteHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "DELETE Command") {
        $deleteHeading = $p
        break
    }
}

if ($deleteHeading -eq $null) {
    throw "Could not find the 'DELETE Command' paragraph"
}

# Create a blank paragraph before it; InsertParagraphBefore() splits the
# paragraph mark so that the *original* $deleteHeading object now refers
# to the newly-created (still empty) paragraph, and the real "DELETE
# Command" text moves to $deleteHeading.Next().
$deleteHeading.Range.InsertParagraphBefore() | Out-Null

$newPara = $deleteHeading

$newParaBody =
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:ind w:right="27"/>' + $rPr + '</w:pPr>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">UPDATE </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>pri</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>marykey</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> IN </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>tablename</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> SET field TO</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> value</w:t></w:r>'

$newParaXml = New-PkgXml ('<w:p>' + $newParaBody + '</w:p>')

$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 2. Split "DELETE primarykey FROM tablename"'s trailing run so that
#    "tablename" is individually wrapped in proofErr markers.
# ---------------------------------------------------------------------

$deleteFromPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "DELETE primarykey FROM tablename") {
        $deleteFromPara = $p
        break
    }
}

if ($deleteFromPara -eq $null) {
    throw "Could not find the 'DELETE primarykey FROM tablename' paragraph"
}

$fromRange = $deleteFromPara.Range
$fromRange.Find.Execute(" FROM tablename") | Out-Null

# Re-seat as a brand-new Range over the same span so InsertXML performs
# a true in-place replace instead of appending after a stale Find range.
$fromRange = $d.Range($fromRange.Start, $fromRange.End)

$fromBody =
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> FROM </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>tablename</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'

$fromXml = New-PkgXml ('<w:p>' + $fromBody + '</w:p>')

$fromRange.InsertXML($fromXml)

Write-Host "UPDATE command bullet inserted and DELETE...FROM tablename proofed."
